$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values in row 38 (year 2023) ---
$ws.Range("B38").Value = 210
$ws.Range("E38").Value = 20
$ws.Range("F38").Value = 620
$ws.Range("G38").Value = 295.162537574381
$ws.Range("H38").Value = 2600
$ws.Range("I38").Value = 720

# --- Update data values in row 40 (year 2025) ---
$ws.Range("C40").Value = 245.384669457003
$ws.Range("J40").Value = 620.379947410462

# --- Remove the conditional formatting rule that was tied to B38 ---
# (the B38 value is no longer treated as a missing/NA placeholder, so the
# "not blank" highlight rule that singled it out is no longer needed)
$cells = $ws.Cells
$fc = $cells.FormatConditions
$totalBefore = $fc.Count

$targetIndex = -1
for ($i = 1; $i -le $fc.Count; $i++) {
    if ($fc.Item($i).AppliesTo.Address() -eq '$B$38') {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $removedRawPriority = $totalBefore + 1 - $targetIndex

    $fc.Item($targetIndex).Delete()

    $newCount = $fc.Count
    for ($i = 1; $i -le $newCount; $i++) {
        if ($i -lt $targetIndex) {
            $origIndex = $i
        } else {
            $origIndex = $i + 1
        }
        $origRaw = $totalBefore + 1 - $origIndex
        if ($origRaw -gt $removedRawPriority) {
            $newRaw = $origRaw - 1
        } else {
            $newRaw = $origRaw
        }
        $fc.Item($i).Priority = $newRaw
    }
}
